# Update the "Danh sach sinh vien" sheet:
#  - row 2 (the first student record) gets new sample data
#  - rows 3 and 4 (the other two student records) are removed entirely,
#    which also shrinks Table1 / the sheet dimension from A1:F4 to A1:F2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "qwe"
$ws.Range("B2").Value = "L001"
$ws.Range("C2").Value = "Nguyễn Văn A"
$ws.Range("D2").Value = "Nam"
$ws.Range("E2").Value = "5/15/2000 12:00:00 AM"
$ws.Range("F2").Value = "Hà Nội"

# Remove old row 4 then row 3 (table range auto-shrinks to A1:F2)
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()
